# Generate Report for Handback
# Adds a new handed-back file ("716fdf9f-d902-4d0a-9342-6582ee004f4b") that is
# in sync with en-US to the Overview / zh-cn / de-de report sheets.

$wb = $excel.ActiveWorkbook

$fileId   = "716fdf9f-d902-4d0a-9342-6582ee004f4b"
$mdName   = "$fileId.md"
$status   = "Handed back: in sync with en-US"
$reason   = "Include"
$hashZh   = "03df20538d4474328b5fb8528df3dfca565da723"

# RGB(100,149,237) == hex 6495ED, stored by Excel's Font.Color as BGR long.
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2   # xlUnderlineStyleSingle
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = $mdName
$ws1.Range("B4").Value = $status
$ws1.Range("C4").Value = $status

$ws1.Hyperlinks.Add(
    $ws1.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8f0a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8/e2e/$mdName",
    "", "", $mdName) | Out-Null
Style-AsHyperlink $ws1.Range("A4")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$fileId.$hashZh.zh-cn.xlf"

$ws2.Range("A4").Value = $mdName
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = $status
$ws2.Range("D4").Value = $zhXlf
$ws2.Range("E4").Value = "2016-03-15 03:27:44"
$ws2.Range("F4").Value = $mdName
$ws2.Range("G4").Value = $zhXlf
$ws2.Range("H4").Value = "2016-03-15 03:28:28"
$ws2.Range("I4").Value = $reason

$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$srcUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e8f0a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8/e2e/$mdName"
$handoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhXlf"
$srcForkUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1/e2e/$mdName"
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhXlf"

$ws2.Hyperlinks.Add($ws2.Range("A4"), $srcUrl, "", "", $mdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), $srcUrl, "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), $handoffUrl, "", "", $zhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F4"), $srcForkUrl, "", "", $mdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G4"), $handbackUrl, "", "", $zhXlf) | Out-Null

Style-AsHyperlink $ws2.Range("A4")
Style-AsHyperlink $ws2.Range("B4")
Style-AsHyperlink $ws2.Range("D4")
Style-AsHyperlink $ws2.Range("F4")
Style-AsHyperlink $ws2.Range("G4")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$hashDe = "03df20538d4474328b5fb8528df3dfca565da723"
$deXlf = "$fileId.$hashDe.de-de.xlf"

$ws3.Range("A4").Value = $mdName
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = $status
$ws3.Range("D4").Value = $deXlf
$ws3.Range("E4").Value = "2016-03-15 03:27:53"
$ws3.Range("F4").Value = $mdName
$ws3.Range("G4").Value = $deXlf
$ws3.Range("H4").Value = "2016-03-15 03:28:42"
$ws3.Range("I4").Value = $reason

$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$srcUrlDe = "https://github.com/OpenLocalizationTest/oltest/blob/e8f0a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8/e2e/$mdName"
$handoffUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deXlf"
$srcForkUrlDe = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2e3f4/e2e/$mdName"
$handbackUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2e3f4a5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deXlf"

$ws3.Hyperlinks.Add($ws3.Range("A4"), $srcUrlDe, "", "", $mdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), $srcUrlDe, "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), $handoffUrlDe, "", "", $deXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F4"), $srcForkUrlDe, "", "", $mdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G4"), $handbackUrlDe, "", "", $deXlf) | Out-Null

Style-AsHyperlink $ws3.Range("A4")
Style-AsHyperlink $ws3.Range("B4")
Style-AsHyperlink $ws3.Range("D4")
Style-AsHyperlink $ws3.Range("F4")
Style-AsHyperlink $ws3.Range("G4")
